$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 299 (pushes the existing 299:318 block down to 300:319),
# matching the "Fruta / hortaliza, semanal" weekly refresh which prepends the
# newest observation ahead of the existing history for this variety block.
$ws.Rows.Item(299).Insert()

$ws.Cells.Item(299, 1).Value = 4
$ws.Cells.Item(299, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(299, 3).Value = "Los Lagos"
$ws.Cells.Item(299, 4).Value = 44585
$ws.Cells.Item(299, 5).Value = 10
$ws.Cells.Item(299, 6).Value = 100114001
$ws.Cells.Item(299, 7).Value = "Papa"
$ws.Cells.Item(299, 8).Value = "Patagonia"
$ws.Cells.Item(299, 9).Value = "1a nueva(o)"
$ws.Cells.Item(299, 10).Value = 300
$ws.Cells.Item(299, 11).Value = 9000
$ws.Cells.Item(299, 12).Value = 10000
$ws.Cells.Item(299, 13).Value = 9500
$ws.Cells.Item(299, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(299, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(299, 16).Value = 380
$ws.Cells.Item(299, 17).Value = 25
$ws.Cells.Item(299, 18).Value = "Hortaliza"
